$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-6: clear the y_0_forecast (C) and y_1_forecast (E) values that were
# erroneously populated too early in the naive forecaster output.
$ws.Range("E2").ClearContents()

$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()

$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Rows 7-19: corrected forecast values from the bugfixed naive forecaster.
$ws.Range("C7").Value = 0.03393100538855442
$ws.Range("E7").Value = 0.3351240474928963

$ws.Range("C8").Value = 1.743169463154315
$ws.Range("E8").Value = 1.13752746419209

$ws.Range("C9").Value = 1.311489985227077
$ws.Range("E9").Value = 1.069982194174801

$ws.Range("C10").Value = 1.721454720714122
$ws.Range("E10").Value = 1.335637690776181

$ws.Range("C11").Value = 1.818507532114921
$ws.Range("E11").Value = 1.41495314213913

$ws.Range("C12").Value = 2.466427116525516
$ws.Range("E12").Value = 1.859723853307749

$ws.Range("C13").Value = 1.128030950601477
$ws.Range("E13").Value = 0.9160236606447159

$ws.Range("C14").Value = 0.1798886261929367
$ws.Range("E14").Value = 0.8029144802146782

$ws.Range("C15").Value = -1.929204335549095
$ws.Range("E15").Value = 0.8768515943972544

$ws.Range("C16").Value = 1.632302710072264
$ws.Range("E16").Value = 0.9471575920676267

$ws.Range("C17").Value = 0.01243672673012508
$ws.Range("E17").Value = 0.5779606211723021

$ws.Range("C18").Value = 0.01966607787367014
$ws.Range("E18").Value = 0.5811853063761419

$ws.Range("C19").Value = 0.5364374648222148
$ws.Range("E19").Value = 0.591050555601802
